$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 69 needs the style of the last existing data row (68) copied for
# the styled columns (A uses style 1, D uses style 2) before values are set.
$ws.Cells.Item(68, 1).Copy()
$ws.Cells.Item(69, 1).PasteSpecial(-4122)
$ws.Cells.Item(68, 4).Copy()
$ws.Cells.Item(69, 4).PasteSpecial(-4122)

$ws.Cells.Item(34, 3).Value = 14000047
$ws.Cells.Item(36, 4).Value = 43656
$ws.Cells.Item(36, 5).Value = "ITCH00001"
$ws.Cells.Item(37, 2).Value = 13000034
$ws.Cells.Item(37, 3).Value = 14000035
$ws.Cells.Item(37, 4).Value = 43656
$ws.Cells.Item(37, 5).Value = "ITCH00002"
$ws.Cells.Item(38, 2).Value = 13000034
$ws.Cells.Item(38, 3).Value = 14000036
$ws.Cells.Item(38, 4).Value = 43656
$ws.Cells.Item(38, 5).Value = "ITCH00003"
$ws.Cells.Item(39, 2).Value = 13000034
$ws.Cells.Item(39, 3).Value = 14000038
$ws.Cells.Item(39, 4).Value = 43656
$ws.Cells.Item(39, 5).Value = "ITCH00004"
$ws.Cells.Item(40, 4).Value = 43657
$ws.Cells.Item(40, 5).Value = "ITCH00005"
$ws.Cells.Item(41, 4).Value = 43657
$ws.Cells.Item(41, 5).Value = "ITCH00006"
$ws.Cells.Item(42, 2).Value = 13000040
$ws.Cells.Item(42, 3).Value = 14000052
$ws.Cells.Item(42, 4).Value = 43767
$ws.Cells.Item(42, 5).Value = "ITCH00016"
$ws.Cells.Item(43, 2).Value = 13000040
$ws.Cells.Item(43, 3).Value = 14000053
$ws.Cells.Item(43, 4).Value = 43767
$ws.Cells.Item(43, 5).Value = "ITCH00017"
$ws.Cells.Item(44, 3).Value = 14000054
$ws.Cells.Item(44, 5).Value = "ITCH00018"
$ws.Cells.Item(45, 3).Value = 14000055
$ws.Cells.Item(45, 5).Value = "ITCH00019"
$ws.Cells.Item(46, 3).Value = 14000056
$ws.Cells.Item(46, 5).Value = "ITCH00020"
$ws.Cells.Item(47, 3).Value = 14000058
$ws.Cells.Item(47, 5).Value = "ITCH00022"
$ws.Cells.Item(48, 3).Value = 14000059
$ws.Cells.Item(48, 5).Value = "ITCH00023"
$ws.Cells.Item(49, 3).Value = 14000060
$ws.Cells.Item(49, 5).Value = "ITCH00024"
$ws.Cells.Item(50, 3).Value = 14000066
$ws.Cells.Item(50, 5).Value = "ITCH00028"
$ws.Cells.Item(51, 2).Value = 13000041
$ws.Cells.Item(51, 3).Value = 14000041
$ws.Cells.Item(51, 4).Value = 43731
$ws.Cells.Item(51, 5).Value = "ITCH00007"
$ws.Cells.Item(52, 2).Value = 13000042
$ws.Cells.Item(52, 3).Value = 14000042
$ws.Cells.Item(52, 4).Value = 43766
$ws.Cells.Item(52, 5).Value = "ITCH00009"
$ws.Cells.Item(53, 2).Value = 13000042
$ws.Cells.Item(53, 4).Value = 43766
$ws.Cells.Item(53, 5).Value = "ITCH00010"
$ws.Cells.Item(54, 2).Value = 13000042
$ws.Cells.Item(54, 3).Value = 14000046
$ws.Cells.Item(54, 4).Value = 43766
$ws.Cells.Item(54, 5).Value = "ITCH00011"
$ws.Cells.Item(55, 2).Value = 13000042
$ws.Cells.Item(55, 3).Value = 14000049
$ws.Cells.Item(55, 4).Value = 43766
$ws.Cells.Item(55, 5).Value = "ITCH00013"
$ws.Cells.Item(56, 2).Value = 13000042
$ws.Cells.Item(56, 3).Value = 14000050
$ws.Cells.Item(56, 5).Value = "ITCH00014"
$ws.Cells.Item(57, 2).Value = 13000042
$ws.Cells.Item(57, 3).Value = 14000051
$ws.Cells.Item(57, 5).Value = "ITCH00015"
$ws.Cells.Item(58, 2).Value = 13000042
$ws.Cells.Item(58, 3).Value = 14000057
$ws.Cells.Item(58, 5).Value = "ITCH00021"
$ws.Cells.Item(59, 2).Value = 13000042
$ws.Cells.Item(59, 3).Value = 14000069
$ws.Cells.Item(59, 5).Value = "ITCH00030"
$ws.Cells.Item(60, 2).Value = 13000061
$ws.Cells.Item(60, 3).Value = 14000061
$ws.Cells.Item(60, 4).Value = 43815
$ws.Cells.Item(60, 5).Value = "ITCH00025"
$ws.Cells.Item(61, 2).Value = 13000063
$ws.Cells.Item(61, 3).Value = 14000063
$ws.Cells.Item(61, 4).Value = 43817
$ws.Cells.Item(61, 5).Value = "ITCH00026"
$ws.Cells.Item(62, 2).Value = 13000063
$ws.Cells.Item(62, 3).Value = 14000068
$ws.Cells.Item(62, 4).Value = 43817
$ws.Cells.Item(62, 5).Value = "ITCH00029"
$ws.Cells.Item(63, 2).Value = 13000063
$ws.Cells.Item(63, 3).Value = 14000070
$ws.Cells.Item(63, 4).Value = 43817
$ws.Cells.Item(63, 5).Value = "ITCH00031"
$ws.Cells.Item(64, 2).Value = 13000064
$ws.Cells.Item(64, 3).Value = 14000064
$ws.Cells.Item(64, 4).Value = 43780
$ws.Cells.Item(64, 5).Value = "ITCH00027"
$ws.Cells.Item(65, 2).Value = 13000071
$ws.Cells.Item(65, 3).Value = 14000071
$ws.Cells.Item(65, 4).Value = 43774
$ws.Cells.Item(65, 5).Value = "ITCH00032"
$ws.Cells.Item(66, 2).Value = 13000072
$ws.Cells.Item(66, 3).Value = 14000072
$ws.Cells.Item(66, 4).Value = 43804
$ws.Cells.Item(66, 5).Value = "ITCH00033"
$ws.Cells.Item(67, 2).Value = 13000073
$ws.Cells.Item(67, 3).Value = 14000073
$ws.Cells.Item(67, 4).Value = 43844
$ws.Cells.Item(67, 5).Value = "ITCH00034"
$ws.Cells.Item(68, 2).Value = 13000073
$ws.Cells.Item(68, 3).Value = 14000074
$ws.Cells.Item(68, 4).Value = 43844
$ws.Cells.Item(68, 5).Value = "ITCH00035"
$ws.Cells.Item(69, 1).Value = 67
$ws.Cells.Item(69, 2).Value = 13000075
$ws.Cells.Item(69, 3).Value = 14000075
$ws.Cells.Item(69, 4).Value = 43845
$ws.Cells.Item(69, 5).Value = "ITCH00036"
